# Updates the LR-pairs (F2 -> F2rl1) sheet with freshly recomputed TPM-based
# NATMI statistics. The "Target cluster" dimension now also includes "MuSCs"
# (previously only ECs/FAPs were present), so every Sending-cluster x
# Target-cluster combination is rewritten, and new rows are appended for the
# combinations that did not exist before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "F2", "F2rl1", "ECs", 3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 2, 0.6666666666666666, 0.4424383333333333, 1.327315, 0.06873529267527652, 0.06873529267527652, 0.731893436835, 6.587040931515, 0.02823361217733215, 0.02823361217733215),
    @("ECs", "F2", "F2rl1", "FAPs", 3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 3, 1, 5.857218, 17.571654, 0.9099518806603505, 0.9099518806603506, 9.689168160485998, 87.20251344437399, 0.3737705551058091, 0.3737705551058092),
    @("ECs", "F2", "F2rl1", "MuSCs", 3, 1, 1.654227, 4.962681, 0.4107585939979205, 0.4107585939979205, 1, 0.3333333333333333, 0.1371873333333333, 0.411562, 0.02131282666437293, 0.02131282666437293, 0.226938990858, 2.042450917722, 0.008754426714779214, 0.008754426714779214),
    @("FAPs", "F2", "F2rl1", "ECs", 3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 2, 0.6666666666666666, 0.4424383333333333, 1.327315, 0.06873529267527652, 0.06873529267527652, 0.5701817835633333, 5.131636052069999, 0.02199540334358256, 0.02199540334358256),
    @("FAPs", "F2", "F2rl1", "FAPs", 3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 3, 1, 5.857218, 17.571654, 0.9099518806603505, 0.9099518806603506, 7.548349124267999, 67.93514211841199, 0.2911860539087374, 0.2911860539087374),
    @("FAPs", "F2", "F2rl1", "MuSCs", 3, 1, 1.288726, 3.866178, 0.3200015957958394, 0.3200015957958394, 1, 0.3333333333333333, 0.1371873333333333, 0.411562, 0.02131282666437293, 0.02131282666437293, 0.1767968833373333, 1.591171950036, 0.006820138543519454, 0.006820138543519454),
    @("MuSCs", "F2", "F2rl1", "ECs", 3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 2, 0.6666666666666666, 0.4424383333333333, 1.327315, 0.06873529267527652, 0.06873529267527652, 0.3919558245411111, 3.52760242087, 0.01512013659884069, 0.01512013659884069),
    @("MuSCs", "F2", "F2rl1", "FAPs", 3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 3, 1, 5.857218, 17.571654, 0.9099518806603505, 0.9099518806603506, 5.188905521388, 46.700149692492, 0.2001678642579684, 0.2001678642579684),
    @("MuSCs", "F2", "F2rl1", "MuSCs", 3, 1, 0.8858993333333333, 2.657698, 0.2199763179924491, 0.2199763179924491, 1, 0.3333333333333333, 0.1371873333333333, 0.411562, 0.02131282666437293, 0.02131282666437293, 0.1215341671417778, 1.093807504276, 0.004688317135640048, 0.004688317135640048),
    @("Resolving-Mac", "F2", "F2rl1", "ECs", 2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 2, 0.6666666666666666, 0.4424383333333333, 1.327315, 0.06873529267527652, 0.06873529267527652, 0.08777814305944444, 0.7900032875349999, 0.003386140555521128, 0.003386140555521128),
    @("Resolving-Mac", "F2", "F2rl1", "FAPs", 2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 3, 1, 5.857218, 17.571654, 0.9099518806603505, 0.9099518806603506, 1.162050574734, 10.458455172606, 0.04482740738783561, 0.04482740738783562),
    @("Resolving-Mac", "F2", "F2rl1", "MuSCs", 2, 0.6666666666666666, 0.1983963333333333, 0.595189, 0.04926349221379096, 0.04926349221379096, 1, 0.3333333333333333, 0.1371873333333333, 0.411562, 0.02131282666437293, 0.02131282666437293, 0.02721746391311111, 0.244957175218, 0.001049944270434212, 0.001049944270434212)
)

$numRows = $data.Count
$numCols = $data[0].Count

for ($i = 0; $i -lt $numRows; $i++) {
    $rowNum = $i + 2
    for ($j = 0; $j -lt $numCols; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $data[$i][$j]
    }
}
